# Append a new paragraph at the very end of the document:
#   "Realidades invertebrales que no tienen sentido."
# "invertebrales" is flagged by the spell checker (like "source", "tree",
# "santia" and "puhsh" elsewhere in this document), so it is wrapped in
# proofErr spellStart/spellEnd markers, split across three runs. We build
# that structure directly as a WordprocessingML fragment and insert it at
# the end of the document's content.

$d = $word.ActiveDocument

$end = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:r><w:t xml:space="preserve">Realidades </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>invertebrales</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> que no tienen sentido.</w:t></w:r>' +
       '</w:p>'

$end.InsertXML($xml)
